# Weekly price update: a new "Albahaca" price-report row for
# Terminal La Palmera de La Serena is inserted at row 27, pushing the
# existing rows 27-36 down to 28-37 (dimension grows from R36 to R37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 27, shifting rows
# 27..36 down to 28..37.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly entry.
$ws.Cells.Item(27, 1).Value  = 8
$ws.Cells.Item(27, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(27, 3).Value  = "Coquimbo"
$ws.Cells.Item(27, 4).Value  = 44460
$ws.Cells.Item(27, 5).Value  = 4
$ws.Cells.Item(27, 6).Value  = 100112052
$ws.Cells.Item(27, 7).Value  = "Albahaca"
$ws.Cells.Item(27, 8).Value  = "Sin especificar"
$ws.Cells.Item(27, 9).Value  = "Primera"
$ws.Cells.Item(27, 10).Value = 800
$ws.Cells.Item(27, 11).Value = 4000
$ws.Cells.Item(27, 12).Value = 4500
$ws.Cells.Item(27, 13).Value = 4250
$ws.Cells.Item(27, 14).Value = "$/paquete"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 4250
$ws.Cells.Item(27, 17).Value = 1
$ws.Cells.Item(27, 18).Value = "Hortaliza"
